# Refresh the crypto price/volume snapshot (cryptos list update,
# Fri Oct 13 09:49:52 UTC 2023 / GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces Excel to store these as text (matching the original
# inline-string cells) instead of silently coercing numeric-looking values
# like "206.68" / "21.50" into Doubles and losing the exact formatting.
$ws.Range("D2").Value = "'26.910.43"
$ws.Range("E2").Value = "'  +0.51%  "
$ws.Range("D3").Value = "'1.554.58"
$ws.Range("E3").Value = "'  +0.22%  "
$ws.Range("E4").Value = "'  +0.34%  "
$ws.Range("D5").Value = "'206.68"
$ws.Range("E5").Value = "'  +1.03%  "
$ws.Range("E6").Value = "'  +0.96%  "
$ws.Range("E7").Value = "'  +0.32%  "
$ws.Range("E8").Value = "'  +0.73%  "
$ws.Range("D9").Value = "'21.50"
$ws.Range("E9").Value = "'  +0.33%  "
$ws.Range("E10").Value = "'  +0.09%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "'  -0.06%  "
$ws.Range("D12").Value = "'1.777.22"
$ws.Range("E12").Value = "'  +0.22%  "
$ws.Range("D13").Value = "'1.555.17"
$ws.Range("E13").Value = "'  +0.03%  "
$ws.Range("E14").Value = "'  +1.06%  "
$ws.Range("D15").Value = "'0.514"
$ws.Range("E15").Value = "'  +0.81%  "
$ws.Range("D16").Value = "'26.908.33"
$ws.Range("E16").Value = "'  +0.52%  "
$ws.Range("E17").Value = "'  +0.88%  "
$ws.Range("D18").Value = "'214.46"
$ws.Range("E18").Value = "'  +0.36%  "
$ws.Range("D19").Value = "'0.0₃0688"
$ws.Range("E19").Value = "'  +1.14%  "
$ws.Range("D20").Value = "'7.24"
$ws.Range("E20").Value = "'  -0.37%  "
$ws.Range("D21").Value = "'1.01"
$ws.Range("E21").Value = "'  +0.31%  "
$ws.Range("E22").Value = "'  -0.70%  "
$ws.Range("E23").Value = "'  +1.39%  "
$ws.Range("E24").Value = "'  -2.03%  "
$ws.Range("D25").Value = "'153.26"
$ws.Range("E25").Value = "'  +0.89%  "
$ws.Range("D26").Value = "'6.66"
$ws.Range("E26").Value = "'  +2.01%  "
$ws.Range("E27").Value = "'  +0.11%  "
$ws.Range("E28").Value = "'  +0.33%  "
$ws.Range("E30").Value = "'  -0.81%  "
$ws.Range("E31").Value = "'  -0.45%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("E32").Value = "'  +2.19%  "
$ws.Range("D33").Value = "'1.368.26"
$ws.Range("E33").Value = "'  +0.28%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("E34").Value = "'  +1.76%  "
$ws.Range("E35").Value = "'  +3.51%  "
$ws.Range("E36").Value = "'  +5.25%  "
$ws.Range("E37").Value = "'  +0.47%  "
$ws.Range("E38").Value = "'  +0.97%  "
$ws.Range("E39").Value = "'  +0.48%  "
$ws.Range("E40").Value = "'  +0.76%  "
$ws.Range("E41").Value = "'  +0.31%  "
$ws.Range("D42").Value = "'0.993"
$ws.Range("E42").Value = "'  +0.50%  "
$ws.Range("E43").Value = "'  -0.41%  "
$ws.Range("D44").Value = "'2.26"
$ws.Range("E44").Value = "'  +3.41%  "
$ws.Range("D45").Value = "'63.64"
$ws.Range("E45").Value = "'  +1.14%  "
$ws.Range("E46").Value = "'  -1.33%  "
$ws.Range("D47").Value = "'1.689.83"
$ws.Range("E47").Value = "'  -0.09%  "
$ws.Range("D48").Value = "'86.08"
$ws.Range("E48").Value = "'  +0.00%  "
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("D50").Value = "'0.0956"
$ws.Range("E50").Value = "'  +1.15%  "
$ws.Range("E51").Value = "'  +0.41%  "
